$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update weekly content text (edit order matches author's editing sequence
# so that new shared-string entries land in the same order)
$ws.Range("E13").Value = "Coursework support session "
$ws.Range("E3").Value = "Objects, Variables & Operators"
$ws.Range("E4").Value = "Control Structures: Conditional Statements"
$ws.Range("E5").Value = "Control Structures: Loops"
$ws.Range("E10").Value = "Reading and plotting data, Matplotlib"
$ws.Range("F11").Value = "CW set"

# Remove the now-unused column H (style-only, no data)
$ws.Columns("H").Delete()

# Hide the Concept/Theme helper columns
$ws.Columns("C").Hidden = $true
$ws.Columns("D").Hidden = $true

# Re-draw the top border on the merged "Programming Fundamentals" block header
$ws.Range("C3:D3").Borders(8).LineStyle = 1
$ws.Range("C3:D3").Borders(8).Weight = -4138

# Match the author's final selection
$ws.Range("F12").Select()
